# chamadaBelaVista.xlsx - fix chamada, testes e ajustes na edição das turmas

$wb = $excel.ActiveWorkbook

# ===================================================================
# Sheet "Alunos": correct row 17 category + add 4 new students (18-21)
# ===================================================================
$alunos = $wb.Worksheets.Item("Alunos")

# Row 17 "Janja da Silva" Categoria was wrongly set, fix it
$alunos.Range("K17").Value = "Não definida"

# Make sure new Data de Nascimento cells use the same date format as
# the rest of the column before writing the date values into them.
$alunos.Range("I18:I21").NumberFormat = $alunos.Range("I2").NumberFormat

# Row 18 - Lula da silva
$alunos.Range("A18").Value = "Lula da silva"
$alunos.Range("B18").Value = "Sim"
$alunos.Range("C18").Value = "Segunda e Quarta"
$alunos.Range("D18").Value = "11h00"
$alunos.Range("E18").Value = "Claudia"
$alunos.Range("F18").Value = "Iniciação"
$alunos.Range("G18").Value = 73
$alunos.Range("H18").Value = "Masculino"
$alunos.Range("I18").Value = "1952-02-12"
$alunos.Range("J18").Value = "(11) 9 9877-1313"
$alunos.Range("K18").Value = "Não definida"
$alunos.Range("L18").Value = "11h00"

# Row 19 - Alexandre Cabeça de Pica
$alunos.Range("A19").Value = "Alexandre Cabeça de Pica"
$alunos.Range("B19").Value = "Sim"
$alunos.Range("C19").Value = "Segunda e Quarta"
$alunos.Range("D19").Value = "11h00"
$alunos.Range("E19").Value = "Claudia"
$alunos.Range("F19").Value = "Iniciação"
$alunos.Range("G19").Value = 52
$alunos.Range("H19").Value = "Masculino"
$alunos.Range("I19").Value = "1973-01-13"
$alunos.Range("J19").Value = "(11) 9 9931-1533"
$alunos.Range("K19").Value = "Não definida"
$alunos.Range("L19").Value = "11h00"

# Row 20 - Dolores da Farmácia
$alunos.Range("A20").Value = "Dolores da Farmácia"
$alunos.Range("B20").Value = "Sim"
$alunos.Range("C20").Value = "Segunda e Quarta"
$alunos.Range("D20").Value = "10h15"
$alunos.Range("E20").Value = "Claudia"
$alunos.Range("F20").Value = "-"
$alunos.Range("G20").Value = 57
$alunos.Range("H20").Value = "Feminino"
$alunos.Range("I20").Value = "1968-01-01"
$alunos.Range("J20").Value = "(19) 9 9961-2344"
$alunos.Range("K20").Value = "Não definida"
$alunos.Range("L20").Value = "10h15"

# Row 21 - Jacir Novais
$alunos.Range("A21").Value = "Jacir Novais"
$alunos.Range("B21").Value = "Sim"
$alunos.Range("C21").Value = "Segunda e Quarta"
$alunos.Range("D21").Value = "10h15"
$alunos.Range("E21").Value = "Claudia"
$alunos.Range("F21").Value = "Nível 1"
$alunos.Range("G21").Value = 48
$alunos.Range("H21").Value = "Masculino"
$alunos.Range("I21").Value = "1977-01-30"
$alunos.Range("J21").Value = "(19) 9 9861-2346"
$alunos.Range("K21").Value = "F45+"
$alunos.Range("L21").Value = "09h30"

# ===================================================================
# Sheet "Turmas": swap Adulto A/B for rows 18-19, add Claudia turmas
# ===================================================================
$turmas = $wb.Worksheets.Item("Turmas")

$turmas.Range("C18").Value = "Adulto B"
$turmas.Range("C19").Value = "Adulto A"

# Row 29
$turmas.Range("A29").Value = "Segunda e Quarta"
$turmas.Range("B29").Value = "10h15"
$turmas.Range("C29").Value = "Nível 1"
$turmas.Range("D29").Value = "Claudia"
$turmas.Range("E29").Value = "10h15"

# Row 30
$turmas.Range("A30").Value = "Segunda e Quarta"
$turmas.Range("B30").Value = "09h30"
$turmas.Range("C30").Value = "Nível 2"
$turmas.Range("D30").Value = "Claudia"
$turmas.Range("E30").Value = "09h30"

# F29/G29/F30/G30 stay empty, matching the (blank) Atalho/Data de
# Início cells already used by every other row (e.g. F28/G28).
$turmas.Range("F28").Copy()
$turmas.Range("F29:G30").PasteSpecial(-4122)

# ===================================================================
# Sheet "Registros": add December attendance columns AA..AT
# ===================================================================
$registros = $wb.Worksheets.Item("Registros")

$novasDatas = @("01/12/2025","03/12/2025","05/12/2025","06/12/2025","07/12/2025","08/12/2025","10/12/2025","12/12/2025","13/12/2025","14/12/2025","15/12/2025","17/12/2025","20/12/2025","21/12/2025","22/12/2025","26/12/2025","27/12/2025","28/12/2025","29/12/2025","31/12/2025")
$novasColunas = @("AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT")

# Force text so the dd/mm/yyyy strings are not auto-converted to real
# Excel dates (matches the literal-text header cells already in row 1).
$headerRange = $registros.Range("AA1:AT1")
$headerRange.NumberFormat = "@"
for ($i = 0; $i -lt $novasColunas.Length; $i++) {
    $registros.Range($novasColunas[$i] + "1").Value = $novasDatas[$i]
}
# Copy the header style (bold, border, alignment) from the existing
# last header cell (Z1) onto the new header cells.
$registros.Range("Z1").Copy()
$headerRange.PasteSpecial(-4122)

# Create the (empty) attendance cells for the existing student rows
# 2-8, matching the blank cells already present for columns H..Z.
$registros.Range("A2").Copy()
$registros.Range("AA2:AT8").PasteSpecial(-4122)
